# Update "想去人数" (number of people interested) counts on the
# "展览" and "全部类型" worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 590
$wsExhibit.Range("F3").Value = 129
$wsExhibit.Range("F4").Value = 31
$wsExhibit.Range("F6").Value = 354
$wsExhibit.Range("F7").Value = 1685
$wsExhibit.Range("F8").Value = 96

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 590
$wsAll.Range("F3").Value = 129
$wsAll.Range("F4").Value = 31
$wsAll.Range("F6").Value = 354
$wsAll.Range("F7").Value = 105
$wsAll.Range("F11").Value = 1685
$wsAll.Range("F12").Value = 96
